# Updated capital structure database
# Apply updated values to rows 2 and 3 of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D" = 0.00328
    "E" = -0.0781
    "G" = 0.6350098619329388
    "H" = 0.6350098619329388
    "I" = 0.3773175542406312
    "J" = 0.3239595162672086
    "K" = 204
    "L" = 0.2011834319526627
    "M" = 149.5
    "N" = 0.05493697864990996
    "O" = 0.7328431372549019
    "P" = 149.5
    "Q" = 0.05493697864990996
    "R" = 0.7328431372549019
    "U" = 168.5
    "V" = 0.06191893580274133
    "W" = 0.1224930947520115
    "X" = 0.07314346024225007
    "Y" = 0.04934963450976144
    "Z" = 0.6252350797575519
    "AA" = 0.2025508539915461
    "AB" = 0.07311194674092447
    "AC" = 0.1294389072506216
    "AD" = 2.11
    "AF" = 2.11
    "AG" = -166.39
    "AH" = 0.0007747639907322069
    "AI" = 0.001166899862294756
    "AJ" = -0.06512558172303524
    "AK" = -0.1014752608693
    "AM" = -11.3
    "AN" = 0.003849662470352125
    "AP" = -0.3035759897828863
    "AQ" = -33.85840707964602
}

foreach ($row in 2, 3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
